# Leave Card update — 5/10/2023 11:17 AM
# Adds a new "TICC" marker, new leave-earning rows for Dec 2022 - Apr 2023,
# a new "2023" year-section header row, and the corresponding SL entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header area -----------------------------------------------------
# F4 (merged F4:G4) gets the "TICC" label.
$ws.Range("F4").Value = "TICC"

# --- Table1 body (rows continue the monthly EARNED entries) ----------
# Row 87 (Dec 2022 entry, already dated/labelled) now earns 1.25.
$ws.Range("C87").Value = 1.25

# Row 90 (next empty slot) becomes the Dec 1, 2022 entry.
$ws.Range("A90").Value = 44896
$ws.Range("C90").Value = 1.25

# Row 91 becomes the "2023" year-section header, styled like the other
# year header rows (e.g. A75 = "2022"). Copy that row's formatting, then
# force the value to be stored as text (quote-prefixed) just like the
# other year headers which live inside date-formatted cells.
$ws.Range("A75").Copy() | Out-Null
$ws.Range("A91").PasteSpecial(-4122) | Out-Null
$ws.Range("A91").Value = "'2023"

# Row 92 = Jan 1, 2023
$ws.Range("A92").Value = 44927
$ws.Range("C92").Value = 1.25

# Row 93 = Feb 1, 2023
$ws.Range("A93").Value = 44958
$ws.Range("C93").Value = 1.25

# Row 94 = Mar 1, 2023
$ws.Range("A94").Value = 44986
$ws.Range("C94").Value = 1.25

# Row 95 = Apr 1, 2023, with a sick-leave particular, an absence-with-pay
# of 2 hours, and a remark referencing the covered dates.
$ws.Range("A95").Value = 45017
$ws.Range("B95").Value = "SL(2-0-0)"
$ws.Range("C95").Value = 1.25
$ws.Range("H95").Value = 2
$ws.Range("K95").Value = "4/19,20/2023"
